# Updated capital structure database
# Refresh computed ratios for the two Nigeria "Healthcare Support Services"
# rows (row 2: industry aggregate, row 3: Union Diagnostic and Clinical
# Services Plc) to reflect the latest inputs.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---------------------------------------------------------------
$ws.Range("D2").Value  = 0.04099999999999999
$ws.Range("G2").Value  = 0.2078167115902965
$ws.Range("H2").Value  = 0.2078167115902965
$ws.Range("I2").Value  = -0.398921832884097
$ws.Range("J2").Value  = -0.398921832884097
$ws.Range("K2").Value  = -1.55
$ws.Range("L2").Value  = -0.4177897574123989
$ws.Range("O2").Value  = -0
$ws.Range("R2").Value  = -0
$ws.Range("U2").Value  = 0.066
$ws.Range("V2").Value  = 0.02283737024221453
$ws.Range("W2").Value  = -0.1220472440944882
$ws.Range("X2").Value  = 0.08519420346152869
$ws.Range("Y2").Value  = -0.2072414475560169
$ws.Range("Z2").Value  = 0.2967050543825976
$ws.Range("AA2").Value = -0.1183621241202815
$ws.Range("AB2").Value = 0.08519420346152869
$ws.Range("AC2").Value = -0.2035563275818102
$ws.Range("AG2").Value = -0.066
$ws.Range("AJ2").Value = -0.02337110481586402
$ws.Range("AK2").Value = -0.006386684730017419
$ws.Range("AN2").Value = -0
$ws.Range("AP2").Value = 0.07719298245614035

# --- Row 3 ---------------------------------------------------------------
$ws.Range("D3").Value  = 0.04099999999999999
$ws.Range("G3").Value  = 0.2078167115902965
$ws.Range("H3").Value  = 0.2078167115902965
$ws.Range("I3").Value  = -0.398921832884097
$ws.Range("J3").Value  = -0.398921832884097
$ws.Range("K3").Value  = -1.55
$ws.Range("L3").Value  = -0.4177897574123989
$ws.Range("O3").Value  = 0
$ws.Range("R3").Value  = 0
$ws.Range("U3").Value  = 0.066
$ws.Range("V3").Value  = 0.02283737024221453
$ws.Range("W3").Value  = -0.1220472440944882
$ws.Range("X3").Value  = 0.08519420346152869
$ws.Range("Y3").Value  = -0.2072414475560169
$ws.Range("Z3").Value  = 0.2967050543825976
$ws.Range("AA3").Value = -0.1183621241202815
$ws.Range("AB3").Value = 0.08519420346152869
$ws.Range("AC3").Value = -0.2035563275818102
$ws.Range("AG3").Value = -0.066
$ws.Range("AJ3").Value = -0.02337110481586402
$ws.Range("AK3").Value = -0.006386684730017419
$ws.Range("AN3").Value = -0
$ws.Range("AP3").Value = 0.07719298245614035
